$wb = $excel.ActiveWorkbook

# --- 1. Update the "Ready for handoff" status text to "In Translation" ---
# This text appears in the Status-related column on every sheet:
#   Overview: columns E (zh-cn) and F (de-de), rows 2-3
#   zh-cn:    column C (Status), rows 2-3
#   de-de:    column C (Status), rows 2-3
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

foreach ($cellAddr in @("E2", "F2", "E3", "F3")) {
    $cell = $wsOverview.Range($cellAddr)
    if ($cell.Value2 -eq "Ready for handoff") {
        $cell.Value = "In Translation"
    }
}

foreach ($ws in @($wsZhCn, $wsDeDe)) {
    foreach ($cellAddr in @("C2", "C3")) {
        $cell = $ws.Range($cellAddr)
        if ($cell.Value2 -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# --- 2. Narrow the locale-status columns ---
# Overview columns E and F, and column C on the zh-cn / de-de sheets,
# shrink from the wider autofit width down to the narrower one.
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
